$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells K1, L1 (bold/centered like the rest of the header row)
$ws.Range("K1").Value = "FollowUpBuffer"
$ws.Range("L1").Value = "FollowUpProb"
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add new data cells K2, L2
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 0.25

# Adjust column widths to match the target layout
# (ColumnWidth values chosen so the resulting stored OOXML width is the closest
# achievable match to the target widths of 16.5703125 / 17.85546875 / 15.140625)
$ws.Columns.Item(10).ColumnWidth = 15.666666666666666
$ws.Columns.Item(11).ColumnWidth = 17
$ws.Columns.Item(12).ColumnWidth = 14.333333333333334

# Update the active selection to the new last header cell
$ws.Range("L1").Select()
